$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '29.679.24'
Set-TextValue $ws.Range("E2") '  -1.42%  '

Set-TextValue $ws.Range("D3") '2.097.83'
Set-TextValue $ws.Range("E3") '  -0.44%  '

Set-TextValue $ws.Range("E4") '  +0.60%  '

Set-TextValue $ws.Range("D5") '343.52'
Set-TextValue $ws.Range("E5") '  -1.96%  '

Set-TextValue $ws.Range("D7") '0.5159'
Set-TextValue $ws.Range("E7") '  -0.11%  '

Set-TextValue $ws.Range("D8") '0.4379'
Set-TextValue $ws.Range("E8") '  -2.55%  '

Set-TextValue $ws.Range("D9") '53.33'
Set-TextValue $ws.Range("E9") '  +0.59%  '

Set-TextValue $ws.Range("D10") '0.09191'
Set-TextValue $ws.Range("E10") '  +2.52%  '

Set-TextValue $ws.Range("E11") '  -0.68%  '

Set-TextValue $ws.Range("D12") '24.53'
Set-TextValue $ws.Range("E12") '  -4.91%  '

Set-TextValue $ws.Range("D13") '2.089.66'
Set-TextValue $ws.Range("E13") '  -0.94%  '

Set-TextValue $ws.Range("D14") '6.758'
Set-TextValue $ws.Range("E14") '  -0.11%  '

Set-TextValue $ws.Range("D15") '8.175'
Set-TextValue $ws.Range("E15") '  +0.33%  '

Set-TextValue $ws.Range("D16") '101.78'
Set-TextValue $ws.Range("E16") '  +2.21%  '

Set-TextValue $ws.Range("D17") '0.00001153'
Set-TextValue $ws.Range("E17") '  +0.07%  '

Set-TextValue $ws.Range("E18") '  +0.56%  '

Set-TextValue $ws.Range("D19") '20.98'
Set-TextValue $ws.Range("E19") '  +1.04%  '

Set-TextValue $ws.Range("D20") '0.06671'
Set-TextValue $ws.Range("E20") '  +0.02%  '

Set-TextValue $ws.Range("D21") '1.008'
Set-TextValue $ws.Range("E21") '  +0.60%  '

Set-TextValue $ws.Range("D22") '6.201'
Set-TextValue $ws.Range("E22") '  -0.68%  '

Set-TextValue $ws.Range("D23") '29.747.29'
Set-TextValue $ws.Range("E23") '  -1.51%  '

Set-TextValue $ws.Range("D24") '12.46'
Set-TextValue $ws.Range("E24") '  -3.83%  '

Set-TextValue $ws.Range("D25") '2.300'
Set-TextValue $ws.Range("E25") '  -2.22%  '

Set-TextValue $ws.Range("D26") '2.340.52'
Set-TextValue $ws.Range("E26") '  -0.78%  '

Set-TextValue $ws.Range("D28") '161.72'
Set-TextValue $ws.Range("E28") '  -0.83%  '

Set-TextValue $ws.Range("D29") '2.500'
Set-TextValue $ws.Range("E29") '  -2.45%  '

Set-TextValue $ws.Range("D30") '133.55'

Set-TextValue $ws.Range("D31") '1.126'
Set-TextValue $ws.Range("E31") '  -5.32%  '

Set-TextValue $ws.Range("D32") '1.666'
Set-TextValue $ws.Range("E32") '  +0.85%  '

Set-TextValue $ws.Range("E33") '  -1.66%  '

Set-TextValue $ws.Range("D34") '6.195'
Set-TextValue $ws.Range("E34") '  -1.30%  '

Set-TextValue $ws.Range("E35") '  -0.52%  '

Set-TextValue $ws.Range("D36") '6.285'
Set-TextValue $ws.Range("E36") '  +5.61%  '

Set-TextValue $ws.Range("D37") '10.39'
Set-TextValue $ws.Range("E37") '  +1.84%  '

Set-TextValue $ws.Range("D38") '0.02574'
Set-TextValue $ws.Range("E38") '  -0.71%  '

Set-TextValue $ws.Range("D39") '0.06681'
Set-TextValue $ws.Range("E39") '  -2.51%  '

Set-TextValue $ws.Range("D40") '0.6998'
Set-TextValue $ws.Range("E40") '  +2.27%  '

Set-TextValue $ws.Range("E41") '  -0.56%  '

Set-TextValue $ws.Range("D42") '1.328'
Set-TextValue $ws.Range("E42") '  +5.80%  '

Set-TextValue $ws.Range("D43") '0.2219'
Set-TextValue $ws.Range("E43") '  -4.23%  '

Set-TextValue $ws.Range("D44") '0.6834'
Set-TextValue $ws.Range("E44") '  +6.17%  '

Set-TextValue $ws.Range("D45") '14.37'
Set-TextValue $ws.Range("E45") '  -0.06%  '

Set-TextValue $ws.Range("D46") '2.310'

Set-TextValue $ws.Range("D47") '3.620'
Set-TextValue $ws.Range("E47") '  -1.47%  '

Set-TextValue $ws.Range("D48") '0.00000000354'
Set-TextValue $ws.Range("E48") '  -1.95%  '

$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue $ws.Range("D49") '1.217'
Set-TextValue $ws.Range("E49") '  -0.68%  '

$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D50") '1.199'
Set-TextValue $ws.Range("E50") '  +2.52%  '

Set-TextValue $ws.Range("D51") '81.28'
Set-TextValue $ws.Range("E51") '  -3.12%  '
